$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.577.64"
$ws.Range("E2").Value = "  -1.69%  "

$ws.Range("D3").Value = "1.879.42"
$ws.Range("E3").Value = "  -1.90%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.19%  "

$ws.Range("D5").Value = "'234.21"
$ws.Range("E5").Value = "  -4.57%  "

$ws.Range("D6").Value = "'0.9993"
$ws.Range("E6").Value = "  +0.07%  "

$ws.Range("D7").Value = "'0.4865"
$ws.Range("E7").Value = "  -1.92%  "

$ws.Range("D8").Value = "'0.2875"
$ws.Range("E8").Value = "  -4.17%  "

$ws.Range("D9").Value = "'0.06635"
$ws.Range("E9").Value = "  -2.47%  "

$ws.Range("D10").Value = "1.874.91"
$ws.Range("E10").Value = "  -2.11%  "

$ws.Range("D11").Value = "'16.70"
$ws.Range("E11").Value = "  -1.95%  "

$ws.Range("D12").Value = "'0.07212"
$ws.Range("E12").Value = "  -1.39%  "

$ws.Range("D13").Value = "'88.26"
$ws.Range("E13").Value = "  -1.73%  "

$ws.Range("D14").Value = "'4.976"
$ws.Range("E14").Value = "  -1.97%  "

$ws.Range("D15").Value = "'0.6606"
$ws.Range("E15").Value = "  -3.16%  "

$ws.Range("D16").Value = "30.533.84"
$ws.Range("E16").Value = "  -1.61%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.000007800"
$ws.Range("E17").Value = "  -2.95%  "

$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'0.9995"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").Value = "'12.92"
$ws.Range("E19").Value = "  -3.06%  "

$ws.Range("D20").Value = "2.119.75"
$ws.Range("E20").Value = "  -1.83%  "

$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.22%  "

$ws.Range("D22").Value = "'4.717"
$ws.Range("E22").Value = "  -3.23%  "

$ws.Range("D23").Value = "'186.93"
$ws.Range("E23").Value = "  +6.01%  "

$ws.Range("D24").Value = "'6.030"
$ws.Range("E24").Value = "  -0.23%  "

$ws.Range("D25").Value = "'9.261"
$ws.Range("E25").Value = "  -0.81%  "

$ws.Range("D26").Value = "'155.33"
$ws.Range("E26").Value = "  +2.26%  "

$ws.Range("D27").Value = "'18.31"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'1.822"
$ws.Range("E28").Value = "  -6.49%  "

$ws.Range("D29").Value = "'1.397"
$ws.Range("E29").Value = "  -1.54%  "

$ws.Range("D30").Value = "'4.235"
$ws.Range("E30").Value = "  -2.83%  "

$ws.Range("D31").Value = "'0.08994"
$ws.Range("E31").Value = "  +0.60%  "

$ws.Range("D32").Value = "'3.911"
$ws.Range("E32").Value = "  -4.01%  "

$ws.Range("D33").Value = "'0.05190"
$ws.Range("E33").Value = "  -1.61%  "

$ws.Range("D34").Value = "'0.7318"
$ws.Range("E34").Value = "  -2.22%  "

$ws.Range("D35").Value = "'1.074"
$ws.Range("E35").Value = "  -6.08%  "

$ws.Range("D36").Value = "'2.688"
$ws.Range("E36").Value = "  +1.25%  "

$ws.Range("D37").Value = "'0.01800"
$ws.Range("E37").Value = "  -7.00%  "

$ws.Range("D38").Value = "'2.649"
$ws.Range("E38").Value = "  -3.39%  "

$ws.Range("D39").Value = "'0.9200"
$ws.Range("E39").Value = "  -2.55%  "

$ws.Range("D40").Value = "'2.035"
$ws.Range("E40").Value = "  -7.95%  "

$ws.Range("D41").Value = "'0.4295"
$ws.Range("E41").Value = "  -1.91%  "

$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "'103.89"
$ws.Range("E42").Value = "  -1.17%  "

$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D43").Value = "'0.9958"
$ws.Range("E43").Value = "  -0.33%  "

$ws.Range("D44").Value = "'5.677"
$ws.Range("E44").Value = "  -4.98%  "

$ws.Range("D45").Value = "'0.1332"
$ws.Range("E45").Value = "  +0.43%  "

$ws.Range("D46").Value = "'7.240"
$ws.Range("E46").Value = "  -7.60%  "

$ws.Range("D47").Value = "'0.05811"
$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("D48").Value = "'8.618"
$ws.Range("E48").Value = "  +0.19%  "

$ws.Range("D49").Value = "'1.409"
$ws.Range("E49").Value = "  +1.25%  "

$ws.Range("D50").Value = "'0.3882"
$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("E51").Value = "  -1.27%  "
